# Docx writer: Use different style for block quotes in notes.
#
# Adds a new paragraph style "Footnote Block Text" (styleId
# "FootnoteBlockText"), based on "Footnote Text" and followed by
# "Footnote Text", mirroring the existing "Block Text" style's
# paragraph formatting (100 twips before/after spacing; 480 twips
# left/right indent; no first-line indent).

$d = $word.ActiveDocument

$style = $d.Styles.Add("FootnoteBlockText", 1)   # wdStyleTypeParagraph
$style.NameLocal = "Footnote Block Text"
$style.BaseStyle = "Footnote Text"
$style.NextParagraphStyle = "Footnote Text"
$style.Priority = 9
$style.UnhideWhenUsed = $true
$style.QuickStyle = $true

# ParagraphFormat distances are expressed in points over COM; the
# target OOXML values are twentieths-of-a-point (twips), so divide by
# 20: 100 -> 5pt, 480 -> 24pt.
$style.ParagraphFormat.SpaceBefore = 5
$style.ParagraphFormat.SpaceAfter = 5
$style.ParagraphFormat.FirstLineIndent = 0
$style.ParagraphFormat.LeftIndent = 24
$style.ParagraphFormat.RightIndent = 24
